# 5.2.1.1b — add a new "2023" data column (column R) that mirrors the
# existing "2022" column (Q) in formatting, then fill in the real 2023
# figures (or the "no data" placeholder "…") for each row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (styles) + placeholder values for rows 3-25 from column Q
# (2022) into the new column R (2023). This picks up the correct style for
# every row in one shot, including the "…" no-data marker already used in
# rows 10-25.
$ws.Range("Q3:Q25").Copy($ws.Range("R3:R25"))

# Header row: 2023
$ws.Range("R4").Value = 2023

# Actual 2023 figures that differ from the copied (2022) values.
$ws.Range("R5").Value = 11357
$ws.Range("R7").Value = 11002
$ws.Range("R8").Value = 355

# Rows 3, 6 and 9 stay blank (section separators), and rows 10-25 keep the
# "…" placeholder copied above, matching the source data.
